$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit swaps the per-record fields of row 2 and row 3 (columns A, B, E,
# F, G, H, I, Q, R) while leaving the shared/common columns (C, D, J, P, S,
# T, U, V, W, Y, Z, AA, AB, AD, AE, AG, AT, AW, AX, AY) untouched.

# Numeric columns can be swapped directly.
$numericCols = @("A", "B", "E", "Q", "R")
foreach ($col in $numericCols) {
    $cell2 = $ws.Range($col + "2")
    $cell3 = $ws.Range($col + "3")
    $tmp = $cell2.Value2
    $cell2.Value = $cell3.Value2
    $cell3.Value = $tmp
}

# Text columns (species name / latin name / author).
$textCols = @("F", "G", "H")
foreach ($col in $textCols) {
    $cell2 = $ws.Range($col + "2")
    $cell3 = $ws.Range($col + "3")
    $tmp = $cell2.Value2
    $cell2.Value = $cell3.Value2
    $cell3.Value = $tmp
}

# Column I ("Antal") is stored as text even though it looks numeric, so
# force text format before assigning to keep it a string value, then
# restore the default "Normal" style so no stray formatting is left behind.
$i2 = $ws.Range("I2").Value2
$i3 = $ws.Range("I3").Value2

$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = [string]$i3
$ws.Range("I2").Style = "Normal"

$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = [string]$i2
$ws.Range("I3").Style = "Normal"
